$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 157, pushing the existing rows 157-196 down to 158-197.
$ws.Rows(157).Insert()

# Populate the newly inserted row 157 with the new weekly record.
# Non date/number columns mirror the template used by every other row in
# this block (same market/product), matching D,J,K,L,M,P per the data update.
$ws.Range("A157").Value = 5
$ws.Range("B157").Value = "Macroferia Regional de Talca"
$ws.Range("C157").Value = "Maule"
$ws.Range("D157").Value = 44508
$ws.Range("E157").Value = 7
$ws.Range("F157").Value = 100112009
$ws.Range("G157").Value = "Acelga"
$ws.Range("H157").Value = "Sin especificar"
$ws.Range("I157").Value = "Primera"
$ws.Range("J157").Value = 500
$ws.Range("K157").Value = 2000
$ws.Range("L157").Value = 2000
$ws.Range("M157").Value = 2000
$ws.Range("N157").Value = "`$/docena de atados (4 kilos)"
$ws.Range("O157").Value = "Región del Maule"
$ws.Range("P157").Value = 500
$ws.Range("Q157").Value = 4
$ws.Range("R157").Value = "Hortaliza"
